# Update design.md diagram: relabel the three callout text boxes from
# (a)/(b)/(c) to (1)/(2)/(3), widen/recenter them over the bundles they
# annotate, and refresh the cached "today" text of the datetimeFigureOut
# placeholder that PowerPoint re-stamps on the slide master + every layout
# whenever the deck is saved.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 1: the three floating "(a)/(b)/(c)" textboxes -> "(1)/(2)/(3)"
#    Reposition/resize + center-align each one, matching the new callouts.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

$callouts = @(
    @{ Name = "TextBox 304"; Text = "(1)"; Left = 6.0;       Width = 120.0 },
    @{ Name = "TextBox 305"; Text = "(2)"; Left = 216.0;     Width = 162.0 },
    @{ Name = "TextBox 306"; Text = "(3)"; Left = 435.8961;  Width = 236.103937 }
)

foreach ($c in $callouts) {
    $sh = $s.Shapes.Item($c.Name)
    $sh.TextFrame.TextRange.Text = $c.Text
    $sh.TextFrame.TextRange.ParagraphFormat.Alignment = 2   # ppAlignCenter
    $sh.Left = $c.Left
    $sh.Width = $c.Width
}

# ---------------------------------------------------------------------
# 2. Refresh the cached datetimeFigureOut field text ("today's date")
#    on the slide master and on every one of its custom layouts.
# ---------------------------------------------------------------------
$newDate = "1/28/2017"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
